$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (even_MAG-GUT13322.fa) entirely - this shifts row 3 (GUT14159) up to row 2
$ws.Rows.Item(2).Delete()

# Now remaining rows are: row1=header, row2=GUT14159 (was row3), row3=GUT23746 (was row4),
# row4=GUT30091 (was row5), row5=GUT54583 (was row6)
# Delete rows 3 through 5 (previously rows 4,5,6) to leave only header + GUT14159 row
$ws.Range("A3:A5").EntireRow.Delete()
